$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.452.66'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '3.894.81'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''601.19'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''167.41'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('D7').Value = '3.893.94'
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.529'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('D11').Value = '''6.48'
$ws.Range('E11').Value = '  +2.84%  '
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('E13').Value = '  +3.71%  '
$ws.Range('D14').Value = '''37.32'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '4.549.46'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '3.872.33'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = '68.538.26'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '''7.46'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '''17.31'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').Value = '''11.05'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').Value = '''490.52'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').Value = '''0.726'
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').Value = '''84.79'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '''2.23'
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').Value = '''11.98'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('D28').Value = '''10.14'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '''2.94'
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').Value = '4.047.99'
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('D32').Value = '''2.36'
$ws.Range('E32').Value = '  -0.77%  '
$ws.Range('D33').Value = '''7.71'
$ws.Range('D34').Value = '''31.72'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '3.852.57'
$ws.Range('E35').Value = '  +2.80%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('D40').Value = '''3.16'
$ws.Range('E40').Value = '  +3.95%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').Value = '''428.52'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '''47.93'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').Value = '''8.56'
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''0.000273'
$ws.Range('E48').Value = '  +19.61%  '
$ws.Range('D49').Value = '''142.76'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').Value = '2.799.74'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').Value = '''39.20'
$ws.Range('E51').Value = '  -0.25%  '
